$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow

# --- Add new log entry row (row 33) ---------------------------------------
$ws.Range("B33").Value = "DiscordiaAgency_Demo_2017_09_23-4.exe"
$ws.Range("C33").Value = "Entwicklung"
$ws.Range("D33").Value = "Anna Franziska"
$ws.Range("E33").Value = "neue Steuerung; Spieler jetzt mit noch weniger Drag"

# --- Update the frozen-pane scroll position / selection --------------------
# Re-establish the freeze at the top row (keeps ySplit=1, state=frozen,
# activePane=bottomLeft) and move the view/selection down near the new row,
# matching the refreshed scroll position after the new entry was appended.
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$win.ScrollRow = 28

$ws.Range("E34").Select()
